{"js": "// Replace the header date and each three-digit \u00f7 one-digit division problem\n// in the table with the new values from the next day's worksheet.\nconst replacements = [\n  [\"2025-04-22 Tuesday\", \"2025-04-23 Wednesday\"],\n  [\"394\u00f73=131, 1\", \"772\u00f79=85, 7\"],\n  [\"878\u00f79=97, 5\", \"676\u00f79=75, 1\"],\n  [\"181\u00f79=20, 1\", \"573\u00f78=71, 5\"],\n  [\"419\u00f78=52, 3\", \"622\u00f74=155, 2\"],\n  [\"543\u00f77=77, 4\", \"442\u00f76=73, 4\"],\n  [\"226\u00f74=56, 2\", \"949\u00f74=237, 1\"],\n  [\"897\u00f75=179, 2\", \"449\u00f72=224, 1\"],\n  [\"994\u00f76=165, 4\", \"161\u00f73=53, 2\"],\n  [\"447\u00f76=74, 3\", \"936\u00f75=187, 1\"],\n  [\"547\u00f75=109, 2\", \"628\u00f77=89, 5\"],\n  [\"195\u00f73=65, 0\", \"296\u00f75=59, 1\"],\n  [\"769\u00f75=153, 4\", \"380\u00f73=126, 2\"],\n  [\"701\u00f78=87, 5\", \"674\u00f72=337, 0\"],\n  [\"367\u00f77=52, 3\", \"178\u00f78=22, 2\"],\n  [\"411\u00f74=102, 3\", \"876\u00f77=125, 1\"],\n  [\"820\u00f77=117, 1\", \"389\u00f73=129, 2\"],\n  [\"581\u00f73=193, 2\", \"778\u00f72=389, 0\"],\n  [\"441\u00f72=220, 1\", \"955\u00f73=318, 1\"],\n  [\"668\u00f72=334, 0\", \"280\u00f72=140, 0\"],\n  [\"343\u00f79=38, 1\", \"657\u00f72=328, 1\"],\n  [\"845\u00f78=105, 5\", \"387\u00f78=48, 3\"],\n  [\"649\u00f79=72, 1\", \"448\u00f72=224, 0\"],\n  [\"978\u00f73=326, 0\", \"225\u00f78=28, 1\"],\n  [\"915\u00f76=152, 3\", \"534\u00f75=106, 4\"],\n  [\"104\u00f73=34, 2\", \"103\u00f79=11, 4\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the header date and each three-digit / one-digit division\n# problem in the table with the new values from the next day's worksheet.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-04-22 Tuesday\", \"2025-04-23 Wednesday\"),\n  @(\"394\u00f73=131, 1\", \"772\u00f79=85, 7\"),\n  @(\"878\u00f79=97, 5\", \"676\u00f79=75, 1\"),\n  @(\"181\u00f79=20, 1\", \"573\u00f78=71, 5\"),\n  @(\"419\u00f78=52, 3\", \"622\u00f74=155, 2\"),\n  @(\"543\u00f77=77, 4\", \"442\u00f76=73, 4\"),\n  @(\"226\u00f74=56, 2\", \"949\u00f74=237, 1\"),\n  @(\"897\u00f75=179, 2\", \"449\u00f72=224, 1\"),\n  @(\"994\u00f76=165, 4\", \"161\u00f73=53, 2\"),\n  @(\"447\u00f76=74, 3\", \"936\u00f75=187, 1\"),\n  @(\"547\u00f75=109, 2\", \"628\u00f77=89, 5\"),\n  @(\"195\u00f73=65, 0\", \"296\u00f75=59, 1\"),\n  @(\"769\u00f75=153, 4\", \"380\u00f73=126, 2\"),\n  @(\"701\u00f78=87, 5\", \"674\u00f72=337, 0\"),\n  @(\"367\u00f77=52, 3\", \"178\u00f78=22, 2\"),\n  @(\"411\u00f74=102, 3\", \"876\u00f77=125, 1\"),\n  @(\"820\u00f77=117, 1\", \"389\u00f73=129, 2\"),\n  @(\"581\u00f73=193, 2\", \"778\u00f72=389, 0\"),\n  @(\"441\u00f72=220, 1\", \"955\u00f73=318, 1\"),\n  @(\"668\u00f72=334, 0\", \"280\u00f72=140, 0\"),\n  @(\"343\u00f79=38, 1\", \"657\u00f72=328, 1\"),\n  @(\"845\u00f78=105, 5\", \"387\u00f78=48, 3\"),\n  @(\"649\u00f79=72, 1\", \"448\u00f72=224, 0\"),\n  @(\"978\u00f73=326, 0\", \"225\u00f78=28, 1\"),\n  @(\"915\u00f76=152, 3\", \"534\u00f75=106, 4\"),\n  @(\"104\u00f73=34, 2\", \"103\u00f79=11, 4\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute($find.Text, $true, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
